$wb = $excel.ActiveWorkbook

function Add-ResultRow {
    param(
        [object]$ws,
        [int]$row,
        [double]$dateValue,
        [double]$h,
        [double]$j,
        [double]$l,
        [double]$n
    )

    $ws.Cells.Item($row, 1).Value = "CSP"
    $ws.Cells.Item($row, 2).Value = $dateValue
    $ws.Cells.Item($row, 2).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($row, 3).Value = 9
    $ws.Cells.Item($row, 4).Value = "Graz"
    $ws.Cells.Item($row, 5).Value = "AT"
    $ws.Cells.Item($row, 6).Value = 0.5
    $ws.Cells.Item($row, 7).Value = 2.5
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = $j
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = $l
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = $n
    $ws.Cells.Item($row, 15).Value = 0
    $ws.Cells.Item($row, 16).Value = "AT"
    $ws.Cells.Item($row, 17).Value = 0.5
    $ws.Cells.Item($row, 18).Value = 2.5
    $ws.Cells.Item($row, 19).Font.Bold = $false
    $ws.Cells.Item($row, 20).Font.Bold = $false
    $ws.Cells.Item($row, 21).Value = 2
    $ws.Cells.Item($row, 22).Value = "Left Hand"
    $ws.Cells.Item($row, 23).Value = "Right Hand"
    $ws.Cells.Item($row, 24).Value = "Feet"
    $ws.Cells.Item($row, 25).Value = "Tongue"
    $ws.Cells.Item($row, 26).Value = 10
}

$wsQDA = $wb.Worksheets.Item("QDA")
Add-ResultRow $wsQDA 8 45412.47420942883 91.666666666666671 73.611111111111114 77.777777777777771 91.666666666666671

$wsKNN = $wb.Worksheets.Item("KNN")
Add-ResultRow $wsKNN 8 45412.474235035661 93.055555555555557 70.833333333333329 69.444444444444443 91.666666666666671

$wsNBPW = $wb.Worksheets.Item("NBPW")
Add-ResultRow $wsNBPW 8 45412.474238313749 97.222222222222229 69.444444444444443 68.055555555555557 88.888888888888886
